# Apply the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.320.15"
$ws.Range("E2").Value = "  +0.00%  "

$ws.Range("D3").Value = "'1.875.67"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'0.7131"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "'242.69"
$ws.Range("E6").Value = "  +0.41%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.08024"
$ws.Range("E8").Value = "  +3.25%  "

$ws.Range("E9").Value = "  +1.61%  "

$ws.Range("D10").Value = "'24.99"
$ws.Range("E10").Value = "  -0.35%  "

$ws.Range("D11").Value = "'0.08227"
$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "'1.893.05"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("B13").Value = "Litecoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D13").Value = "'94.81"
$ws.Range("E13").Value = "  +4.02%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.246"
$ws.Range("E14").Value = "  +0.03%  "

$ws.Range("D15").Value = "'0.7124"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").Value = "'6.404"
$ws.Range("E16").Value = "  +5.64%  "

$ws.Range("D17").Value = "'0.000008551"
$ws.Range("E17").Value = "  +4.41%  "

$ws.Range("D18").Value = "'29.339.92"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").Value = "'243.82"
$ws.Range("E19").Value = "  +1.75%  "

$ws.Range("D20").Value = "'2.137.23"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'7.774"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").Value = "'0.1561"
$ws.Range("E25").Value = "  -1.48%  "

$ws.Range("D26").Value = "'9.048"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").Value = "'162.52"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("E28").Value = "  +0.19%  "

$ws.Range("D29").Value = "'1.501"
$ws.Range("E29").Value = "  -0.44%  "

$ws.Range("D30").Value = "'4.418"
$ws.Range("E30").Value = "  +0.36%  "

$ws.Range("D31").Value = "'4.306"
$ws.Range("E31").Value = "  -0.45%  "

$ws.Range("D32").Value = "'0.05376"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("D33").Value = "'1.167"
$ws.Range("E33").Value = "  -9.24%  "

$ws.Range("E34").Value = "  +0.21%  "

$ws.Range("D35").Value = "'0.7638"
$ws.Range("E35").Value = "  +2.62%  "

$ws.Range("D36").Value = "'1.179"
$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").Value = "'2.688"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("D38").Value = "'0.01877"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").Value = "'1.256.06"
$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("D40").Value = "'2.754"
$ws.Range("E40").Value = "  +0.99%  "

$ws.Range("D41").Value = "'6.477"
$ws.Range("E41").Value = "  +0.57%  "

$ws.Range("D42").Value = "'0.9157"
$ws.Range("E42").Value = "  +3.42%  "

$ws.Range("D43").Value = "'112.84"
$ws.Range("E43").Value = "  +2.84%  "

$ws.Range("D44").Value = "'74.18"
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("E45").Value = "  +9.30%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'2.038.07"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").Value = "'0.5223"
$ws.Range("E48").Value = "  +0.51%  "

$ws.Range("D49").Value = "'1.800"
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").Value = "'9.466"
$ws.Range("E50").Value = "  +0.97%  "

$ws.Range("E51").Value = "  +1.09%  "
